$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 20834400
$ws.Range("I100").Value = 1085.1578
$ws.Range("J100").Value = 100001000
$ws.Range("K100").Value = 1085.1578
$ws.Range("L100").Value = 100001000
$ws.Range("M100").Value = -544.1578
$ws.Range("N100").Value = -100002082
$ws.Range("H105").Value = 32024
$ws.Range("J105").Value = 32024
$ws.Range("L105").Value = 32024
$ws.Range("N105").Value = -39012
$ws.Range("H120").Value = 34666.668
$ws.Range("J120").Value = 34666.668
$ws.Range("L120").Value = 34666.668
$ws.Range("N120").Value = -44342.668
$ws.Range("H129").Value = 1193.4767
$ws.Range("I129").Value = 610.125
$ws.Range("J129").Value = 1253.3077
$ws.Range("K129").Value = 1830.375
$ws.Range("L129").Value = 3759.9231
$ws.Range("M129").Value = 3169.625
$ws.Range("N129").Value = -13759.9231
$ws.Range("H137").Value = 3751453.8
$ws.Range("I137").Value = 1725305.1
$ws.Range("J137").Value = 9093118
$ws.Range("K137").Value = 5175915.300000001
$ws.Range("L137").Value = 27279354
$ws.Range("M137").Value = -5173365.300000001
$ws.Range("N137").Value = -27284454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 30800
$ws.Range("J119").Value = 30800
$ws.Range("L119").Value = 30800
$ws.Range("N119").Value = -40476
$ws.Range("H132").Value = 85035.95
$ws.Range("I132").Value = 93583.5
$ws.Range("J132").Value = 8108
$ws.Range("K132").Value = 280750.5
$ws.Range("L132").Value = 24324
$ws.Range("M132").Value = -278220.5
$ws.Range("N132").Value = -29384

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 15752.5
$ws.Range("I82").Value = 5333.3335
$ws.Range("K82").Value = 5333.3335
$ws.Range("M82").Value = -4950.3335
$ws.Range("H85").Value = 15752.5
$ws.Range("I85").Value = 5333.3335
$ws.Range("K85").Value = 5333.3335
$ws.Range("M85").Value = -4007.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 8725.143
$ws.Range("J51").Value = 11097.2
$ws.Range("L51").Value = 11097.2
$ws.Range("N51").Value = -12569.2
$ws.Range("H61").Value = 8725.143
$ws.Range("J61").Value = 11097.2
$ws.Range("L61").Value = 11097.2
$ws.Range("N61").Value = -11793.2
$ws.Range("H99").Value = 1432.2
$ws.Range("I99").Value = 1315.3572
$ws.Range("K99").Value = 1315.3572
$ws.Range("M99").Value = 182.6428000000001
$ws.Range("H109").Value = 11000
$ws.Range("J109").Value = 11000
$ws.Range("L109").Value = 11000
$ws.Range("N109").Value = -13080
$ws.Range("H110").Value = 35795
$ws.Range("J110").Value = 35795
$ws.Range("L110").Value = 35795
$ws.Range("N110").Value = -43975
$ws.Range("H111").Value = 30000
$ws.Range("J111").Value = 30000
$ws.Range("L111").Value = 30000
$ws.Range("N111").Value = -38180
$ws.Range("H126").Value = 1432.2
$ws.Range("I126").Value = 1315.3572
$ws.Range("K126").Value = 3946.0716
$ws.Range("M126").Value = -1476.0716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 430710.12
$ws.Range("I5").Value = 451.43478
$ws.Range("J5").Value = 1667703.9
$ws.Range("K5").Value = 1354.30434
$ws.Range("L5").Value = 5003111.699999999
$ws.Range("M5").Value = -1242.30434
$ws.Range("N5").Value = -5003335.699999999
$ws.Range("H70").Value = 7188.846
$ws.Range("J70").Value = 7827
$ws.Range("L70").Value = 23481
$ws.Range("N70").Value = -24111
$ws.Range("H73").Value = 7188.846
$ws.Range("J73").Value = 7827
$ws.Range("L73").Value = 23481
$ws.Range("N73").Value = -25665
$ws.Range("H93").Value = 2175.6667
$ws.Range("I93").Value = 500
$ws.Range("J93").Value = 3013.5
$ws.Range("K93").Value = 1500
$ws.Range("L93").Value = 9040.5
$ws.Range("M93").Value = 372
$ws.Range("N93").Value = -12784.5
$ws.Range("H97").Value = 866.6667
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 866.6667
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2600.0001
$ws.Range("M97").Value = ""
$ws.Range("N97").Value = -3592.0001
$ws.Range("H113").Value = 592.75757
$ws.Range("I113").Value = 520.4091
$ws.Range("J113").Value = 737.4545
$ws.Range("K113").Value = 1561.2273
$ws.Range("L113").Value = 2212.3635
$ws.Range("M113").Value = 608.7727
$ws.Range("N113").Value = -6552.3635
$ws.Range("H122").Value = 43527.312
$ws.Range("J122").Value = 47327.582
$ws.Range("L122").Value = 425948.238
$ws.Range("N122").Value = -430848.238
$ws.Range("H135").Value = 430710.12
$ws.Range("I135").Value = 451.43478
$ws.Range("J135").Value = 1667703.9
$ws.Range("K135").Value = 4062.91302
$ws.Range("L135").Value = 15009335.1
$ws.Range("M135").Value = -1527.91302
$ws.Range("N135").Value = -15014405.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1453.3158
$ws.Range("I113").Value = 1328.5714
$ws.Range("J113").Value = 1802.6
$ws.Range("K113").Value = 1328.5714
$ws.Range("L113").Value = 1802.6
$ws.Range("M113").Value = 841.4286
$ws.Range("N113").Value = -6142.6
$ws.Range("H121").Value = 22000
$ws.Range("J121").Value = 22000
$ws.Range("L121").Value = 22000
$ws.Range("N121").Value = -25494
$ws.Range("H123").Value = 29967.75
$ws.Range("J123").Value = 29967.75
$ws.Range("L123").Value = 29967.75
$ws.Range("N123").Value = -34867.75
$ws.Range("H132").Value = 3437.4736
$ws.Range("I132").Value = 3300
$ws.Range("J132").Value = 3735.3333
$ws.Range("K132").Value = 9900
$ws.Range("L132").Value = 11205.9999
$ws.Range("M132").Value = -7370
$ws.Range("N132").Value = -16265.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H132").Value = 1403.5625
$ws.Range("I132").Value = 1124.7742
$ws.Range("J132").Value = 1911.9412
$ws.Range("K132").Value = 3374.3226
$ws.Range("L132").Value = 5735.8236
$ws.Range("M132").Value = -844.3226000000004
$ws.Range("N132").Value = -10795.8236
$ws.Range("H133").Value = 42910.445
$ws.Range("J133").Value = 42910.445
$ws.Range("L133").Value = 42910.445
$ws.Range("N133").Value = -47970.445
